# The project moved into a "funcionalidades_menu" folder layout and two new
# windows were wired up (ventana-notificaciones / VentanaAgendasController).
# On the "notificaciones" sheet of the data workbook, a "Fecha y hora" column
# was introduced between "Titulo" and "Emisor", and the old "Mensaje" column
# was pushed out to the right (after "Emisor"), giving the final header row:
#   Titulo | Fecha y hora | Emisor | Mensaje | Receptor
# The workbook was also left with the "notificaciones" sheet active/selected
# (it used to be "docentes").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("notificaciones")

# Remember the two values that need to move before we overwrite anything.
$mensaje = $ws.Range("B1").Value()
$emisor  = $ws.Range("D1").Value()

# B1: Mensaje -> new "Fecha y hora" header
$ws.Range("B1").Value = "Fecha y hora"
# C1: Fecha -> Emisor (moved left)
$ws.Range("C1").Value = $emisor
# D1: Emisor -> Mensaje (moved right)
$ws.Range("D1").Value = $mensaje
# A1 (Titulo) and E1 (Receptor) stay as they were.

# New "Fecha y hora" / moved "Emisor" column read a bit narrower than the
# previous headers, so give B:C a best-fit-ish width like the resaved file.
$ws.Range("B1:C1").ColumnWidth = 11.220779220779221

# The notificaciones sheet becomes the active tab of the workbook.
$ws.Activate()
